{"js": "// Replace the worksheet date and every three-digit-by-one-digit\n// multiplication answer with the next day's generated values.\n// Each \"before\" string occurs exactly once in the document, so a\n// simple search + full-text replace per pair is safe and unambiguous.\n\nconst replacements = [\n  [\"2024-06-06 Thursday\", \"2024-06-07 Friday\"],\n  [\"307\u00d73=921\", \"243\u00d72=486\"],\n  [\"279\u00d78=2232\", \"443\u00d79=3987\"],\n  [\"915\u00d76=5490\", \"821\u00d78=6568\"],\n  [\"267\u00d77=1869\", \"956\u00d75=4780\"],\n  [\"364\u00d78=2912\", \"495\u00d78=3960\"],\n  [\"943\u00d79=8487\", \"231\u00d74=924\"],\n  [\"757\u00d73=2271\", \"888\u00d75=4440\"],\n  [\"651\u00d73=1953\", \"739\u00d74=2956\"],\n  [\"803\u00d72=1606\", \"498\u00d77=3486\"],\n  [\"176\u00d77=1232\", \"826\u00d78=6608\"],\n  [\"855\u00d77=5985\", \"898\u00d76=5388\"],\n  [\"648\u00d78=5184\", \"678\u00d73=2034\"],\n  [\"795\u00d77=5565\", \"664\u00d74=2656\"],\n  [\"275\u00d74=1100\", \"881\u00d79=7929\"],\n  [\"840\u00d72=1680\", \"765\u00d74=3060\"],\n  [\"525\u00d72=1050\", \"139\u00d72=278\"],\n  [\"436\u00d78=3488\", \"360\u00d78=2880\"],\n  [\"990\u00d79=8910\", \"475\u00d72=950\"],\n  [\"822\u00d76=4932\", \"847\u00d73=2541\"],\n  [\"611\u00d72=1222\", \"318\u00d76=1908\"],\n  [\"222\u00d75=1110\", \"438\u00d78=3504\"],\n  [\"796\u00d76=4776\", \"631\u00d76=3786\"],\n  [\"804\u00d75=4020\", \"332\u00d77=2324\"],\n  [\"803\u00d76=4818\", \"198\u00d72=396\"],\n  [\"836\u00d74=3344\", \"158\u00d74=632\"],\n];\n\nfor (const [before, after] of replacements) {\n  const results = context.document.body.search(before, { matchCase: true });\n  results.load(\"items\");\n  await context.sync();\n\n  for (const range of results.items) {\n    range.insertText(after, \"Replace\");\n  }\n  await context.sync();\n}\n", "ps1": "# Replace the worksheet date and every three-digit-by-one-digit\n# multiplication answer with the next day's generated values.\n# Each \"before\" string occurs exactly once in the document, so a\n# plain Find/Replace (no wildcards) per pair is unambiguous.\n\n$d = $word.ActiveDocument\n\n$replacements = @(\n    @(\"2024-06-06 Thursday\", \"2024-06-07 Friday\"),\n    @(\"307\u00d73=921\", \"243\u00d72=486\"),\n    @(\"279\u00d78=2232\", \"443\u00d79=3987\"),\n    @(\"915\u00d76=5490\", \"821\u00d78=6568\"),\n    @(\"267\u00d77=1869\", \"956\u00d75=4780\"),\n    @(\"364\u00d78=2912\", \"495\u00d78=3960\"),\n    @(\"943\u00d79=8487\", \"231\u00d74=924\"),\n    @(\"757\u00d73=2271\", \"888\u00d75=4440\"),\n    @(\"651\u00d73=1953\", \"739\u00d74=2956\"),\n    @(\"803\u00d72=1606\", \"498\u00d77=3486\"),\n    @(\"176\u00d77=1232\", \"826\u00d78=6608\"),\n    @(\"855\u00d77=5985\", \"898\u00d76=5388\"),\n    @(\"648\u00d78=5184\", \"678\u00d73=2034\"),\n    @(\"795\u00d77=5565\", \"664\u00d74=2656\"),\n    @(\"275\u00d74=1100\", \"881\u00d79=7929\"),\n    @(\"840\u00d72=1680\", \"765\u00d74=3060\"),\n    @(\"525\u00d72=1050\", \"139\u00d72=278\"),\n    @(\"436\u00d78=3488\", \"360\u00d78=2880\"),\n    @(\"990\u00d79=8910\", \"475\u00d72=950\"),\n    @(\"822\u00d76=4932\", \"847\u00d73=2541\"),\n    @(\"611\u00d72=1222\", \"318\u00d76=1908\"),\n    @(\"222\u00d75=1110\", \"438\u00d78=3504\"),\n    @(\"796\u00d76=4776\", \"631\u00d76=3786\"),\n    @(\"804\u00d75=4020\", \"332\u00d77=2324\"),\n    @(\"803\u00d76=4818\", \"198\u00d72=396\"),\n    @(\"836\u00d74=3344\", \"158\u00d74=632\")\n)\n\nforeach ($pair in $replacements) {\n    $findText = $pair[0]\n    $replaceText = $pair[1]\n\n    $range = $d.Content\n    $range.Find.ClearFormatting()\n    $range.Find.Replacement.ClearFormatting()\n    $range.Find.Text = $findText\n    $range.Find.Replacement.Text = $replaceText\n    $range.Find.Forward = $true\n    $range.Find.Wrap = 1\n    $range.Find.MatchCase = $true\n    $range.Find.MatchWholeWord = $false\n    $range.Find.MatchWildcards = $false\n    $range.Find.Execute($findText, $true, $false, $false, $false, $false, $true, 1, $false, $replaceText, 2)\n}\n"}
